# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 currently has the "date only" number format reserved for the most
# recent day (it's the last row in the sheet). Capture that format so it can
# be moved onto the new last row, then restyle row 26 to match the regular
# datetime format used by every earlier row.
$lastRowFormat = $ws.Range("A26").NumberFormat
$ws.Range("A26").NumberFormat = $ws.Range("A25").NumberFormat

# Append the new day's data as row 27.
$ws.Range("A27").Value = 45611
$ws.Range("B27").Value = 68
$ws.Range("C27").Value = 57
$ws.Range("D27").Value = 63

# The newly appended row is now the most recent day, so it gets the
# "date only" format previously used by row 26.
$ws.Range("A27").NumberFormat = $lastRowFormat
